$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report gained a new shortage-item row ("اولويز ماكس طويل جدا") that must be
# inserted as item #33, right before the existing "حبايه" row (old row 39).
# This pushes the three existing trailing item rows down by one, appends the
# totals row and the footer row one row further down, bumps the grand total by
# the new row's sell price, and refreshes the generation timestamp.

# 1) Insert a new blank row at row 39 - this shifts rows 39:44 down to 40:45
#    (data rows, totals row and footer row all move down together), carrying
#    their formatting/merges with them.
$ws.Rows("39:39").Insert()

# 2) Copy formatting (styles/borders/number formats) from the row immediately
#    below (the former row 39, now row 40) onto the freshly inserted row 39 so
#    it matches the other item rows exactly.
$ws.Range("A40:Q40").Copy()
$ws.Range("A39:Q39").PasteSpecial(-4122)

# 3) Fill in the new item's data in row 39.
$ws.Range("A39").Value2 = 33
$ws.Range("C39").Value2 = "اولويز ماكس طويل جدا"
$ws.Range("H39").Value2 = "36:0"
$ws.Range("L39").Value2 = "0"
$ws.Range("N39").Value2 = "35.00"
$ws.Range("P39").Value2 = "35.0000"
$ws.Range("Q39").Value2 = "1:0"

# 4) The item-number column keeps counting sequentially; the old rows already
#    shifted down (now rows 40-42) retain their original numbers (34, 35) and
#    the row that used to be the totals row is now row 43, footer row 44 -
#    update the new last item row (42) to item #36.
$ws.Range("A42").Value2 = 36

# 5) Bump the grand-total cell (now row 43) by the new item's sell price.
$total = $ws.Range("P43").Value2
$ws.Range("P43").Value2 = $total + 35
$ws.Range("P43").NumberFormat = "#.00"

# 6) Refresh the generated-report timestamp shown in the footer (now row 44).
$ws.Range("A44").Value2 = "Saturday, 14 June, 2025 1:54 PM"
